# "Add files via upload"
#
# The sheet "2026" tracks daily injection-molding defect counts. Rows
# 264-285 already hold the 2026-02-06 (serial 46059) block. This edit adds
# the next day's block, 2026-02-18 (serial 46071), into the previously
# blank rows 286-307 - same part list/order as the template block, with
# that day's actual D (input qty) and G:P (defect-by-type) figures.
#
# Columns: B=date, C=part name, D=input qty, E=D-F (good qty),
#          F=SUM(G:P) (total defects), G..P = defect counts by type.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2026")

$dateSerial = 46071   # 2026-02-18

# Per-row part name (shared-string text), input qty (D) and defect
# columns (G..P) that are non-zero. F/E are (re)derived below.
$rows = @(
    @{ Row = 286; Name = "5750 SUBSTRATE TMX 2RB RH";                              D = 264;  Defects = @{} },
    @{ Row = 287; Name = "5749 SUBSTRATE TMX 2RB LH";                              D = 600;  Defects = @{} },
    @{ Row = 288; Name = "WK GRAB HANDLE L/RH OUTER - 1ST(BLACK)";                 D = 808;  Defects = @{ G = 13; H = 5 } },
    @{ Row = 289; Name = "WK GRAB HANDLE L/RH INNER - 2ST";                        D = 808;  Defects = @{ G = 13; H = 5 } },
    @{ Row = 290; Name = "6252 JL HDL RH HANDLE 신규금형";                          D = 3695; Defects = @{ H = 21; L = 18 } },
    @{ Row = 291; Name = "OV1 BRKT-TAIL GATE GARNISH CTR  87395-X9000";            D = 1068; Defects = @{ H = 18 } },
    @{ Row = 292; Name = "NOZZLE-DEFROSTER NO.2(LHD+) / WK";                       D = 1017; Defects = @{ H = 17 } },
    @{ Row = 293; Name = "#6987 TRACER DVR OB LHD";                                D = 2448; Defects = @{} },
    @{ Row = 294; Name = "NQ5PE GRILLECTR SPEAKER";                                D = 420;  Defects = @{} },
    @{ Row = 295; Name = "QV RETAINER T.G.S(RHD)";                                 D = 1001; Defects = @{ G = 1 } },
    @{ Row = 296; Name = "QV BEZEL TGS";                                          D = 312;  Defects = @{} },
    @{ Row = 297; Name = "6243 JL HDL FRT LH HSG w/SWITCH";                        D = 1908; Defects = @{ G = 10; H = 10 } },
    @{ Row = 298; Name = "OV1 BRKT-TAIL GATE GARNISH, LH/RH  87333/4-X9000";       D = 590;  Defects = @{} },
    @{ Row = 299; Name = "#6898-2 WL74 Speaker Grille Subwoofer RH Chrome Ring";   D = 60;   Defects = @{} },
    @{ Row = 300; Name = "BRKT-TAIL GATE GARNISH LH/RH   87733/4-BS000";          D = 1780; Defects = @{} },
    @{ Row = 301; Name = "BRKT RR BEAM UPR,CTR   86632-CH100";                    D = 1000; Defects = @{} },
    @{ Row = 302; Name = "BRKT-RR CORNER RADAR MTG,RH   866G2-BS010";             D = 1200; Defects = @{} },
    @{ Row = 303; Name = "NQ5-PE COVER-RR BUMPER LWR-P1510";                      D = 252;  Defects = @{ P = 2 } },
    @{ Row = 304; Name = "#6922 10.1 in. Display Bezel";                          D = 814;  Defects = @{ G = 5; H = 24; J = 5 } },
    @{ Row = 305; Name = "WD OUTER RING";                                        D = 1622; Defects = @{ H = 2 } },
    @{ Row = 306; Name = "7126 WL Speaker Cover RH,LH";                          D = 720;  Defects = @{} },
    @{ Row = 307; Name = "BRKT-FR BUMPER SIDE UPR,L/RH 86525/6-X9000";           D = 786;  Defects = @{} }
)

# F cells that were overtyped with a literal total instead of keeping the
# =SUM(G:P) formula (breaks the shared-formula chain at those two rows,
# same as in the source file).
$literalF = @{ 288 = 600; 305 = 300 }

# A single non-defect-count formula cell that deviates from a plain value
# (kept as an actual formula, just like in the workbook).
$literalFormulas = @{ "290,J" = "=44-8" }

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 2).Value = $dateSerial          # B - date
    $ws.Cells.Item($row, 3).Value = $r.Name               # C - part name
    $ws.Cells.Item($row, 4).Value = $r.D                  # D - input qty
    $ws.Cells.Item($row, 5).Formula = "=D$row-F$row"      # E - good qty

    if ($literalF.ContainsKey($row)) {
        $ws.Cells.Item($row, 6).Value = $literalF[$row]   # F - literal override
    } else {
        $ws.Cells.Item($row, 6).Formula = "=SUM(G$row`:P$row)"   # F - total defects
    }

    foreach ($col in @("G", "H", "I", "J", "K", "L", "M", "N", "O", "P")) {
        $key = "$row,$col"
        if ($literalFormulas.ContainsKey($key)) {
            $ws.Range("$col$row").Formula = $literalFormulas[$key]
        } elseif ($r.Defects.ContainsKey($col)) {
            $ws.Range("$col$row").Value = $r.Defects[$col]
        }
    }
}

# The workbook was left with the cursor on F306 after the last entry.
$ws.Range("F306").Select()
